$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# --- 1. Wrap the "bbbbbbbbbbbbbbbb" paragraph's run in proofErr spellStart/
#        gramStart ... spellEnd/gramEnd markers (simulates Word's automatic
#        spell/grammar-check bookmarks around the typed word/sentence).
$pB = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*bbbbbbbbbbbbbbbb*") {
        $pB = $p
    }
}
$xmlB = '<w:p xmlns:w="' + $wNs + '">' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:lastRenderedPageBreak/><w:t>bbbbbbbbbbbbbbbb</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '</w:p>'
[void]$pB.Range.InsertXML($xmlB)

# --- 2. Append a new paragraph holding a manual page break, right after the
#        "bbbb" paragraph.
$lastPara = $d.Paragraphs.Last
$rAfterB = $lastPara.Range
$rAfterB.Collapse(0)
[void]$rAfterB.InsertParagraphAfter()
$pBreak = $d.Paragraphs.Last
$xmlBreak = '<w:p xmlns:w="' + $wNs + '"><w:r><w:br w:type="page"/></w:r></w:p>'
[void]$pBreak.Range.InsertXML($xmlBreak)

# --- 3. Append another new paragraph after the page break holding the
#        "ccccc..." text (with its own lastRenderedPageBreak, matching the
#        pattern already used for the "bbbb" paragraph).
$rAfterBreak = $pBreak.Range
$rAfterBreak.Collapse(0)
[void]$rAfterBreak.InsertParagraphAfter()
$pC = $d.Paragraphs.Last
$xmlC = '<w:p xmlns:w="' + $wNs + '"><w:r><w:lastRenderedPageBreak/>' + `
  '<w:t>ccccccccccccccccccccccccccccccccccccccccccccccccccccccccccc</w:t></w:r></w:p>'
[void]$pC.Range.InsertXML($xmlC)

# --- 4. Bump the cached page-number field result in the footer from 2 to 3
#        (document now spans one more page after the edits above).
$footer = $d.Sections(1).Footers.Item(1)
$pageField = $footer.Range.Fields.Item(1)
$pageField.Result.Text = "3"
